{"js": "// Find the paragraph containing the target sentence and split its single\n// run into three runs: the first sentence (struck through), a plain space,\n// and the second sentence (struck through).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst firstSentence =\n  \"Then we can make a bigger array and test multiple of each of the aircraft in the mix.\";\nconst secondSentence =\n  \"The random generation of some amount of each plane can be done at the end.\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(firstSentence) !== -1 && text.indexOf(secondSentence) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst firstResults = target.search(firstSentence, { matchCase: true });\nfirstResults.load(\"items\");\nconst secondResults = target.search(secondSentence, { matchCase: true });\nsecondResults.load(\"items\");\nawait context.sync();\n\nif (firstResults.items.length === 0 || secondResults.items.length === 0) {\n  throw new Error(\"Target sentences not found in paragraph\");\n}\n\nfirstResults.items[0].font.strikeThrough = true;\nsecondResults.items[0].font.strikeThrough = true;\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to the two sentences of the target\n# paragraph, leaving the single space between them unformatted. This\n# naturally splits the original single run into three runs.\n\n$d = $word.ActiveDocument\n\n$firstSentence = \"Then we can make a bigger array and test multiple of each of the aircraft in the mix.\"\n$secondSentence = \"The random generation of some amount of each plane can be done at the end.\"\n\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$found1 = $rng1.Find.Execute($firstSentence)\nif ($found1) {\n    $rng1.Font.StrikeThrough = 1\n}\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute($secondSentence)\nif ($found2) {\n    $rng2.Font.StrikeThrough = 1\n}\n"}
